$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = -22.0112
$ws.Range("A18").Value = -22.16220000000001
$ws.Range("A20").Value = -20.0931
$ws.Range("A27").Value = -22.10099999999999
$ws.Range("A69").Value = -21.63159999999999
$ws.Range("A76").Value = -19.71749999999999
$ws.Range("A82").Value = -21.6257
